# Adds a bit of i18n — set Sheet1 A1:A3 labels (Russian) and shift the
# numeric values into column B, and select Sheet1 as the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Move existing values to column B, add label strings to column A.
$ws1.Range("B1").Value = 0.5
$ws1.Range("B2").Value = 8
$ws1.Range("B3").Value = 1

$ws1.Range("A1").Value = "Коэффициент скорости обучения"
$ws1.Range("A2").Value = "Количество входов нейронной сети"
$ws1.Range("A3").Value = "Размерность выходного слоя"

# Make Sheet1 the active / selected sheet with A1:A3 selected.
$ws1.Activate()
$ws1.Range("A1:A3").Select()

$wb.Save()
